# Append two new sections to the end of the GWAAP project report:
#   - Heading2 "Engineering and TDD concerns in Django" + 4 body paragraphs
#   - Heading2 "Convention over configuration" + 2 body paragraphs
#
# The new material is inserted immediately after the document's final
# existing paragraph ("...committed the new stable version."). While
# typing that material, Word's automatic "_GoBack" (last edit location)
# bookmark would end up re-anchored inside the freshly typed text -
# specifically splitting the word "project" into "proje" / "ct" partway
# through the first new body paragraph - so we reproduce that exactly:
# the old bookmark is removed from the tail of the last paragraph and
# an equivalent one is woven into the new text at that same spot.

$d = $word.ActiveDocument

# Find the end of the document's current last paragraph of real text.
$anchor = $d.Content
$anchor.Find.ClearFormatting()
$found = $anchor.Find.Execute(
    "committed the new stable version.", $false, $false, $false,
    $false, $false, $true, 1, $false, $null, 0)
if (-not $found) {
    throw "Could not locate insertion anchor text"
}
$insertPos = $anchor.End

# That old "_GoBack" bookmark currently sits right at the end of this
# paragraph; drop it here so it can be rebuilt further down, inside
# the new text, in the same operation.
$d.Bookmarks("_GoBack").Delete()

# Build a fresh Range at that exact spot (a brand-new Range object
# avoids a stale/late-bound Range quirk after the Find+Delete above)
# and inject the new paragraphs as real OOXML, so paragraph styles,
# proofing marks (spelling/grammar squiggle anchors), the mid-page
# render-break marker, and the relocated bookmark all land precisely
# where the original author's edit left them.
$insertRange = $d.Range($insertPos, $insertPos)

$newContentXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/new.xml" pkg:contentType="text/xml">
    <pkg:xmlData>
      <w:body xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
          <w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>Engineering and TDD concerns in Django</w:t></w:r></w:p>
          <w:p><w:r><w:t>One of the most difficult aspects of forcing Django development into a software process is that big-picture decisions made during architecture and even iteration planning were often incompatible with Django implementation.  The clearest example of this from the GWAAP proje</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>ct was the split between User and Applicant when referring to the two classes of users who needed access to the system.  From an architectural standpoint, these classes were outlined as two separate subclasses of a joint parent class.  The two distinct classes share some similar needs – the ability to be authenticated by the system, the need for a method of defining permissions and grouping, and the conceptual similarity that both classes represent people who need to use various aspects of the GWAAP system.  For this reason, they needed to derive from the same base class.  However, much of their functionality was also divergent: Applicants need Applicant objects and need to be jailed into only accessing information directly related to them and their application, while Users have much broader access to the system to say nothing of the fact that including an Application object for a faculty User makes no sense.  Thus, during the architecture phase of the project, the two classes were given a hierarchy in keeping with the above.</w:t></w:r></w:p>
          <w:p><w:r><w:t>During implementation, however, the reality of the way Django prefers to handle such things made it impossible to abide by the original architecture.  Django includes a built-in User class (an unfortunate conflict with the GWAAP User class name).  Naming conflicts were the least of the problems, however.  Django “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>automagically</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">” provides Manager </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>objects</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> go with all Django models, of which Django User is a </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">distant subclass.  Although this was resolved from a syntactic standpoint using special Python importing language (“import…from…as”), it made it impossible for the compiler to recognize special Django User methods when operating on GWAAP users.  Specifically, it made it impossible to access the special </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>User.create_</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>user</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">) method, meaning all GWAAP Users and Applicants have to revert to the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Manager.create</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>() method, a less-powerful method that meant many user-specific details had to be set with multiple statements and an explicit save() call.</w:t></w:r></w:p>
          <w:p><w:r><w:t>Worse still was the problem of checking authentication differences between GWAAP User and Applicant objects.  While the Python interpreter recognizes the difference in the two classes, the Django framework does not.  In practice, this means that Django treats all instances of GWAAP User and GWAAP Applicant as interchangeable with Django User.  Thus, the built-in Django methods of authentication (like method decorators) gave the same access to Applicant as to User.  This was an obvious semantic problem and required custom authentication checks for each View.</w:t></w:r></w:p>
          <w:p><w:r><w:t xml:space="preserve">Ideally, the architecture of the system could have been rewritten to accommodate this limitation.  However, at the point in construction this problem manifested itself, there were tens of test cases that depended upon the User/Applicant split.  Changing that assumption would have necessitated not only a rewrite of production code but also an almost total rewrite of test code – in effect, it would have meant starting over from scratch.  Do to the realities of time </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>pressure,</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> I elected to maintain the User/Applicant split from a test-case perspective, but to thereafter handle authentication differences by the inclusion of additional Django Permissions.  This was yet another place in which Django automations were supplanted by hand-crafted code.</w:t></w:r></w:p>
          <w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>Convention over configuration</w:t></w:r></w:p>
          <w:p><w:r><w:t xml:space="preserve">An oft-cited tenet of the Ruby language is the idea of “convention over configuration”.  The idea behind many of Ruby’s design decisions is to make minimal boilerplate code that accomplishes all standard behavior, and leave systems to define only behavior in which they differ from convention.  [PERHAPS USE THE MAC METAPHOR HERE INSTEAD?]  </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Django’s</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> approach is similar in some ways.  If what you want to accomplish is part of the built-in functionality of Django, it is simple, powerful, and secure to do so.  However, designs that run counter to assumptions the framework makes are penalized heavily.  They are difficult to implement, difficult to debug, and difficult to integrate with Django-friendly aspects of the system.</w:t></w:r></w:p>
          <w:p><w:r><w:t>Obviously, this is difficult to reconcile with traditional method of architecture and design in software process.  Although you may construct an accurate, well-planned design from a “green field” perspective, it breaks down as you try to reconcile your design with the preconceptions of Django.</w:t></w:r></w:p>
      </w:body>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$insertRange.InsertXML($newContentXml) | Out-Null
